$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add a new "2022-Q4" sheet, placed right after "总计" and before the
#    existing "2021-Q4" sheet. We create it by duplicating the current
#    "2021-Q4" sheet (so it inherits the exact same layout/formatting),
#    then rename the duplicate and overwrite its data with the new
#    2022-Q4 fund figures. The original "2021-Q4" sheet (and the sheets
#    after it) keep their data untouched, they simply shift right by one
#    tab position.
# ---------------------------------------------------------------------
$wsOldQ4 = $wb.Worksheets.Item("2021-Q4")
$wsOldQ4.Copy($wsOldQ4)

$wsNewQ4 = $wb.Worksheets.Item(2)
$wsNewQ4.Name = "2022-Q4"

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $wsNewQ4.Range("B2") "001468"
$wsNewQ4.Range("C2").Value = "广发改革先锋灵活配置混合"
Set-TextValue $wsNewQ4.Range("D2") "5.72"
Set-TextValue $wsNewQ4.Range("E2") "93.29"
Set-TextValue $wsNewQ4.Range("F2") "2.60"
Set-TextValue $wsNewQ4.Range("G2") "0.1487"

# ---------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: append a new row for 2020-Q4 at
#    the bottom, shift the existing quarterly figures down one row, and
#    write the new 2022-Q4 totals into the top data row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("A4").Copy($wsTotal.Range("A5"))
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = $wsTotal.Range("B4").Value()
$wsTotal.Range("C5").Value = $wsTotal.Range("C4").Value()
$wsTotal.Range("D5").Value = $wsTotal.Range("D4").Value()

$wsTotal.Range("B4").Value = $wsTotal.Range("B3").Value()
$wsTotal.Range("C4").Value = $wsTotal.Range("C3").Value()
$wsTotal.Range("D4").Value = $wsTotal.Range("D3").Value()

$wsTotal.Range("B3").Value = $wsTotal.Range("B2").Value()
$wsTotal.Range("C3").Value = $wsTotal.Range("C2").Value()
$wsTotal.Range("D3").Value = $wsTotal.Range("D2").Value()

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.15
